$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 16:22"

# --- Province/city reordering + updated case numbers ---
# Tenerife jumps ahead of Cordoba and Burgos (new, higher case count)
$ws.Range("A32").Value = "Tenerife"
$ws.Range("B32").Value = 1080
$ws.Range("C32").Value = 157
$ws.Range("D32").Value = 1649
$ws.Range("E32").Value = 56

$ws.Range("A33").Value = "Cordoba"
$ws.Range("B33").Value = 1046
$ws.Range("C33").Value = 84
$ws.Range("D33").Value = 925
$ws.Range("E33").Value = 37

$ws.Range("A34").Value = "Burgos"
$ws.Range("B34").Value = 1003
$ws.Range("C34").Value = 360
$ws.Range("D34").Value = 530
$ws.Range("E34").Value = 113

# Huesca moves ahead of Gran Canaria (Gran Canaria corrected downward)
$ws.Range("A47").Value = "Huesca"
$ws.Range("B47").Value = 417
$ws.Range("C47").Value = 57
$ws.Range("D47").Value = 323
$ws.Range("E47").Value = 37

$ws.Range("A48").Value = "Gran Canaria"
$ws.Range("B48").Value = 411
$ws.Range("C48").Value = 157
$ws.Range("D48").Value = 1649
$ws.Range("E48").Value = 24

# La Palma moves ahead of Lanzarote
$ws.Range("A56").Value = "La Palma"
$ws.Range("B56").Value = 66
$ws.Range("C56").Value = 157
$ws.Range("D56").Value = 1649
$ws.Range("E56").Value = 3

$ws.Range("A57").Value = "Lanzarote"
$ws.Range("B57").Value = 58
$ws.Range("C57").Value = 157
$ws.Range("D57").Value = 1649
$ws.Range("E57").Value = 2

# --- Simple numeric updates (no reordering) ---
$ws.Range("B59").Value = 23
$ws.Range("B62").Value = 7
$ws.Range("B64").Value = 1
$ws.Range("D64").Value = 1407
